$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the rows/columns that need to move
# (D = Fecha, M = Volumen, N = Precio minimo, O = Precio maximo,
#  P = Precio promedio ponderado, R = Origen, S = Precio $/Kg)
$cols = @("D", "M", "N", "O", "P", "R", "S")

$orig = @{}
foreach ($row in 2, 3, 4, 6) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$row").Value2
    }
    $orig[$row] = $rowVals
}

# Target row gets the values previously held by source row
# (a 4-cycle: 2 <- 6 <- 3 <- 4 <- 2)
$mapping = @{ 2 = 6; 3 = 4; 4 = 2; 6 = 3 }

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    foreach ($col in $cols) {
        $ws.Range("$col$target").Value2 = $orig[$source][$col]
    }
}
